# Auto-generated edit script: update cryptocurrency price/volume table
# Applies the cell-level text changes described by the commit diff
# (price updates, volume % updates, and a couple of row-content swaps).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value, and whether the value
# must be forced to Text (using a leading apostrophe / quote-prefix)
# so Excel does not silently reinterpret it as a Number and drop
# significant trailing zeros / alter its representation.
$updates = @(
    @{ Cell = "D2"; Value = '69.900.78'; ForceText = $false }
    @{ Cell = "E2"; Value = '  -1.14%  '; ForceText = $false }
    @{ Cell = "D3"; Value = '3.571.73'; ForceText = $false }
    @{ Cell = "E3"; Value = '  -2.20%  '; ForceText = $false }
    @{ Cell = "E4"; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = "D5"; Value = '575.15'; ForceText = $true }
    @{ Cell = "E5"; Value = '  -3.09%  '; ForceText = $false }
    @{ Cell = "D6"; Value = '186.99'; ForceText = $true }
    @{ Cell = "E6"; Value = '  -3.82%  '; ForceText = $false }
    @{ Cell = "D7"; Value = '3.565.50'; ForceText = $false }
    @{ Cell = "E7"; Value = '  -2.21%  '; ForceText = $false }
    @{ Cell = "D8"; Value = '0.622'; ForceText = $true }
    @{ Cell = "E8"; Value = '  -3.96%  '; ForceText = $false }
    @{ Cell = "E9"; Value = '  +0.03%  '; ForceText = $false }
    @{ Cell = "D10"; Value = '0.184'; ForceText = $true }
    @{ Cell = "E10"; Value = '  +2.12%  '; ForceText = $false }
    @{ Cell = "D11"; Value = '0.650'; ForceText = $true }
    @{ Cell = "E11"; Value = '  -3.63%  '; ForceText = $false }
    @{ Cell = "D12"; Value = '54.98'; ForceText = $true }
    @{ Cell = "E12"; Value = '  -5.86%  '; ForceText = $false }
    @{ Cell = "D13"; Value = '0.0000304'; ForceText = $true }
    @{ Cell = "E13"; Value = '  +3.59%  '; ForceText = $false }
    @{ Cell = "D14"; Value = '9.56'; ForceText = $true }
    @{ Cell = "E14"; Value = '  -3.92%  '; ForceText = $false }
    @{ Cell = "D15"; Value = '4.142.73'; ForceText = $false }
    @{ Cell = "E15"; Value = '  -2.05%  '; ForceText = $false }
    @{ Cell = "D16"; Value = '19.62'; ForceText = $true }
    @{ Cell = "E16"; Value = '  -1.64%  '; ForceText = $false }
    @{ Cell = "D17"; Value = '3.567.05'; ForceText = $false }
    @{ Cell = "E17"; Value = '  -2.10%  '; ForceText = $false }
    @{ Cell = "D18"; Value = '69.833.66'; ForceText = $false }
    @{ Cell = "E18"; Value = '  -1.15%  '; ForceText = $false }
    @{ Cell = "D19"; Value = '12.52'; ForceText = $true }
    @{ Cell = "E19"; Value = '  -2.15%  '; ForceText = $false }
    @{ Cell = "E20"; Value = '  -0.52%  '; ForceText = $false }
    @{ Cell = "D21"; Value = '1.03'; ForceText = $true }
    @{ Cell = "E21"; Value = '  -3.63%  '; ForceText = $false }
    @{ Cell = "D22"; Value = '486.07'; ForceText = $true }
    @{ Cell = "E22"; Value = '  -0.83%  '; ForceText = $false }
    @{ Cell = "D23"; Value = '19.76'; ForceText = $true }
    @{ Cell = "E23"; Value = '  +3.22%  '; ForceText = $false }
    @{ Cell = "D24"; Value = '4.89'; ForceText = $true }
    @{ Cell = "E24"; Value = '  -8.13%  '; ForceText = $false }
    @{ Cell = "D25"; Value = '95.55'; ForceText = $true }
    @{ Cell = "E25"; Value = '  +4.80%  '; ForceText = $false }
    @{ Cell = "D26"; Value = '4.35'; ForceText = $true }
    @{ Cell = "E26"; Value = '  -3.71%  '; ForceText = $false }
    @{ Cell = "B27"; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = "C27"; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false }
    @{ Cell = "D27"; Value = '11.16'; ForceText = $true }
    @{ Cell = "E27"; Value = '  -4.11%  '; ForceText = $false }
    @{ Cell = "B28"; Value = 'ImmutableX'; ForceText = $false }
    @{ Cell = "C28"; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false }
    @{ Cell = "D28"; Value = '2.96'; ForceText = $true }
    @{ Cell = "E28"; Value = '  -7.25%  '; ForceText = $false }
    @{ Cell = "D29"; Value = '9.31'; ForceText = $true }
    @{ Cell = "E29"; Value = '  -3.46%  '; ForceText = $false }
    @{ Cell = "D30"; Value = '31.75'; ForceText = $true }
    @{ Cell = "E30"; Value = '  -3.46%  '; ForceText = $false }
    @{ Cell = "D31"; Value = '7.50'; ForceText = $true }
    @{ Cell = "E31"; Value = '  -5.64%  '; ForceText = $false }
    @{ Cell = "D32"; Value = '67.20'; ForceText = $true }
    @{ Cell = "E32"; Value = '  +2.22%  '; ForceText = $false }
    @{ Cell = "D33"; Value = '12.05'; ForceText = $true }
    @{ Cell = "E33"; Value = '  -1.91%  '; ForceText = $false }
    @{ Cell = "E34"; Value = '  -5.23%  '; ForceText = $false }
    @{ Cell = "D35"; Value = '568.16'; ForceText = $true }
    @{ Cell = "E35"; Value = '  -9.85%  '; ForceText = $false }
    @{ Cell = "B36"; Value = 'Fetch.AI'; ForceText = $false }
    @{ Cell = "C36"; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; ForceText = $false }
    @{ Cell = "D36"; Value = '3.14'; ForceText = $true }
    @{ Cell = "E36"; Value = '  +9.34%  '; ForceText = $false }
    @{ Cell = "B37"; Value = 'InjectiveProtocol'; ForceText = $false }
    @{ Cell = "C37"; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; ForceText = $false }
    @{ Cell = "D37"; Value = '38.43'; ForceText = $true }
    @{ Cell = "E37"; Value = '  -5.61%  '; ForceText = $false }
    @{ Cell = "B38"; Value = 'Dai'; ForceText = $false }
    @{ Cell = "C38"; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; ForceText = $false }
    @{ Cell = "D38"; Value = '0.999'; ForceText = $true }
    @{ Cell = "E38"; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = "D39"; Value = '0.0₃0798'; ForceText = $false }
    @{ Cell = "E39"; Value = '  -3.26%  '; ForceText = $false }
    @{ Cell = "D40"; Value = '0.393'; ForceText = $true }
    @{ Cell = "E40"; Value = '  -5.04%  '; ForceText = $false }
    @{ Cell = "D41"; Value = '3.26'; ForceText = $true }
    @{ Cell = "E41"; Value = '  +11.46%  '; ForceText = $false }
    @{ Cell = "D42"; Value = '3.52'; ForceText = $true }
    @{ Cell = "E42"; Value = '  -1.82%  '; ForceText = $false }
    @{ Cell = "E43"; Value = '  -8.19%  '; ForceText = $false }
    @{ Cell = "D44"; Value = '3.268.02'; ForceText = $false }
    @{ Cell = "E44"; Value = '  -1.17%  '; ForceText = $false }
    @{ Cell = "D45"; Value = '3.00'; ForceText = $true }
    @{ Cell = "E45"; Value = '  -5.50%  '; ForceText = $false }
    @{ Cell = "E46"; Value = '  +4.28%  '; ForceText = $false }
    @{ Cell = "D47"; Value = '0.0439'; ForceText = $true }
    @{ Cell = "E47"; Value = '  -3.43%  '; ForceText = $false }
    @{ Cell = "D48"; Value = '9.59'; ForceText = $true }
    @{ Cell = "E48"; Value = '  +3.67%  '; ForceText = $false }
    @{ Cell = "D49"; Value = '0.136'; ForceText = $true }
    @{ Cell = "E49"; Value = '  -2.15%  '; ForceText = $false }
    @{ Cell = "D50"; Value = '0.998'; ForceText = $true }
    @{ Cell = "E50"; Value = '  -0.07%  '; ForceText = $false }
    @{ Cell = "D51"; Value = '3.18'; ForceText = $true }
    @{ Cell = "E51"; Value = '  -4.52%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe forces Excel to store the value as text
        # (quote-prefixed), preserving the exact digits/trailing zeros.
        $range.Value = "'" + $u.Value
    } else {
        $range.Value = $u.Value
    }
}
